# capitalization fixes for NonGeoLookupTables
# - Fix header capitalization "province_Code" -> "province_code" on ResourceType_TBS
# - Add new Ontario ("ON") rows to Format_TBS and ResourceType_TBS lookup tables
# - Leave Update_TBS / SubjectNTopic_TBS data untouched, adjust active-sheet selection

$wb = $excel.ActiveWorkbook

$wsFormat       = $wb.Worksheets.Item("Format_TBS")
$wsUpdate       = $wb.Worksheets.Item("Update_TBS")
$wsResourceType = $wb.Worksheets.Item("ResourceType_TBS")
$wsSubjectTopic = $wb.Worksheets.Item("SubjectNTopic_TBS")

# ---------------------------------------------------------------------------
# Format_TBS: append new "ON" rows (47-58)
# ---------------------------------------------------------------------------
$formatRows = @(
    @("other", "other", "ON"),
    @("pdf", "PDF", "ON"),
    @("txt", "TXT", "ON"),
    @("xml", "XML", "ON"),
    @("zip", "ZIP", "ON"),
    @("xlsx", "XLSX", "ON"),
    @("docx", "DOCX", "ON"),
    @("xls", "XLS", "ON"),
    @("csv", "CSV", "ON"),
    @("html", "HTML", "ON"),
    @("doc", "DOC", "ON"),
    @("application/msaccess", "other", "ON")
)

$startRow = 47
for ($i = 0; $i -lt $formatRows.Count; $i++) {
    $r = $startRow + $i
    $row = $formatRows[$i]
    $wsFormat.Cells.Item($r, 1).Value = $row[0]
    $wsFormat.Cells.Item($r, 2).Value = $row[1]
    $wsFormat.Cells.Item($r, 3).Value = $row[2]
}

# ---------------------------------------------------------------------------
# ResourceType_TBS: fix header capitalization + append new "ON" rows (7-12)
# ---------------------------------------------------------------------------
$wsResourceType.Range("C1").Value = "province_code"

$resourceTypeRows = @(
    @("unknown", "dataset", "ON"),
    @("application", "application", "ON"),
    @("data", "dataset", "ON"),
    @("technical_document", "dataset", "ON"),
    @("data_dictionary", "guide", "ON"),
    @("information", "guide", "ON")
)

$startRow = 7
for ($i = 0; $i -lt $resourceTypeRows.Count; $i++) {
    $r = $startRow + $i
    $row = $resourceTypeRows[$i]
    $wsResourceType.Cells.Item($r, 1).Value = $row[0]
    $wsResourceType.Cells.Item($r, 2).Value = $row[1]
    $wsResourceType.Cells.Item($r, 3).Value = $row[2]
}

# Column A of ResourceType_TBS needs to widen now that it holds longer values
# (matches the width Excel settles on after auto-fitting the "application/msaccess" entry)
$wsResourceType.Columns.Item(1).ColumnWidth = 23.666666666666664

# ---------------------------------------------------------------------------
# Update view / selection state
# ---------------------------------------------------------------------------
$wsFormat.Range("A59").Select() | Out-Null

$wsUpdate.Application.ActiveWindow.ScrollRow = 16
$wsUpdate.Range("B21").Select() | Out-Null

$wsSubjectTopic.Range("A15").Select() | Out-Null

# ResourceType_TBS becomes the active/selected sheet
$wsResourceType.Activate()
$wsResourceType.Range("C1").Select() | Out-Null

$wb.Save()
